$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows for Score / Best labels first so shared strings are appended in proper order
$ws.Range("A75").Value = "BASIC_TEXT_SCORE"
$ws.Range("B75").Value = "Score"
$ws.Range("C75").Value = "XXXX"
$ws.Range("D75").Value = "XXXX"
$ws.Range("E75").Value = "XXXX"

$ws.Range("A76").Value = "BASIC_TEXT_BEST"
$ws.Range("B76").Value = "Best"
$ws.Range("C76").Value = "XXXX"
$ws.Range("D76").Value = "XXXX"
$ws.Range("E76").Value = "XXXX"

# Update Game Over message (B69) and Game Over body (B70)
$ws.Range("B69").Value = "Game Over"
$ws.Range("B70").Value = "You Survived {0} Turns*2n*Citizen happiness fallen below 20%, citizens no longer feel safe under your control"

$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
$ws.Range("B70").Select()
